$wb = $excel.ActiveWorkbook

# --- Update the "Logs" sheet: append a new row (17) with the new mail log entry ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A17").Value = "Interne taak"
$logs.Range("B17").Value = "kwaliteit@testbedrijf123.nl"
$logs.Range("C17").Value = "Leg dit even neer bij Koen."
$logs.Range("D17").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E17").Value = "Bedankt, we hebben dit doorgestuurd naar support@testbedrijf123.nl."
$logs.Range("F17").Value = "2025-08-14 20:57:06"
$logs.Range("G17").Value = "Nee"
$logs.Range("H17").Value = "Ja"
$logs.Range("I17").Value = "Nee"
$logs.Range("J17").Value = "Nee"

# --- Extend the conditional formatting ranges on "Logs" so they cover the new row too ---
$ranges = @("D2:D16", "G2:G16", "H2:H16", "I2:I16", "J2:J16")
$newRanges = @("D2:D17", "G2:G17", "H2:H17", "I2:I17", "J2:J17")

for ($r = 0; $r -lt $ranges.Length; $r++) {
    $fcs = $logs.Range($ranges[$r]).FormatConditions
    $target = $logs.Range($newRanges[$r])
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($target)
    }
}

# --- Update the "Dashboard" sheet: bump the count for the corresponding category ---
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 11
